$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 223 ("Vega Monumental Concepción" /
# Naranja / Valencia block) to add this week's price readings. This shifts all
# existing data rows 223-268 down to 225-270.
$ws.Rows.Item(223).Insert()
$ws.Rows.Item(223).Insert()

# New row 223: Naranja - New Hall - Primera
$ws.Range("A223").Value = 11
$ws.Range("B223").Value = "Vega Monumental Concepción"
$ws.Range("C223").Value = "Bíobío"
$ws.Range("D223").Value = 44694
$ws.Range("E223").Value = 8
$ws.Range("F223").Value = "Fruta"
$ws.Range("G223").Value = 100102
$ws.Range("H223").Value = "Cítricos"
$ws.Range("I223").Value = 100102005
$ws.Range("J223").Value = "Naranja"
$ws.Range("K223").Value = "New Hall"
$ws.Range("L223").Value = "Primera"
$ws.Range("M223").Value = 220
$ws.Range("N223").Value = 8500
$ws.Range("O223").Value = 9000
$ws.Range("P223").Value = 8773
$ws.Range("Q223").Value = "$/caja 15 kilos granel"
$ws.Range("R223").Value = "Región de O'Higgins"
$ws.Range("S223").Value = 585
$ws.Range("T223").Value = 15

# New row 224: Naranja - Valencia - Primera
$ws.Range("A224").Value = 11
$ws.Range("B224").Value = "Vega Monumental Concepción"
$ws.Range("C224").Value = "Bíobío"
$ws.Range("D224").Value = 44694
$ws.Range("E224").Value = 8
$ws.Range("F224").Value = "Fruta"
$ws.Range("G224").Value = 100102
$ws.Range("H224").Value = "Cítricos"
$ws.Range("I224").Value = 100102005
$ws.Range("J224").Value = "Naranja"
$ws.Range("K224").Value = "Valencia"
$ws.Range("L224").Value = "Primera"
$ws.Range("M224").Value = 180
$ws.Range("N224").Value = 9000
$ws.Range("O224").Value = 10000
$ws.Range("P224").Value = 9444
$ws.Range("Q224").Value = "$/caja 15 kilos granel"
$ws.Range("R224").Value = "Región de O'Higgins"
$ws.Range("S224").Value = 630
$ws.Range("T224").Value = 15
